$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.320.97'
$ws.Range("E2").Value = '  -0.10%  '

$ws.Range("D3").Value = '''1.843.55'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").Value = '''0.9984'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '''240.07'
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("D6").Value = '''0.6285'
$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").Value = '''1.000'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '''0.07452'
$ws.Range("E8").Value = '  -1.89%  '

$ws.Range("D9").Value = '''0.2895'
$ws.Range("E9").Value = '  -0.36%  '

$ws.Range("D10").Value = '''24.39'
$ws.Range("E10").Value = '  -1.00%  '

$ws.Range("D11").Value = '''0.07738'
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").Value = '''1.843.73'
$ws.Range("E12").Value = '  -2.36%  '

$ws.Range("D13").Value = '''4.980'
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("D14").Value = '''0.6791'

$ws.Range("D15").Value = '''0.00001042'
$ws.Range("E15").Value = '  -0.75%  '

$ws.Range("E16").Value = '  -1.44%  '

$ws.Range("D17").Value = '''6.187'
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("D18").Value = '''29.373.80'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("D19").Value = '''227.66'
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("D20").Value = '''12.30'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").Value = '''0.9998'
$ws.Range("E21").Value = '  +0.14%  '

$ws.Range("D22").Value = '''7.513'
$ws.Range("E22").Value = '  +0.72%  '

$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").Value = '''159.12'
$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("D25").Value = '''8.472'
$ws.Range("E25").Value = '  +0.48%  '

$ws.Range("D27").Value = '''17.49'
$ws.Range("E27").Value = '  -1.04%  '

$ws.Range("D28").Value = '''0.06499'
$ws.Range("E28").Value = '  +16.13%  '

$ws.Range("D29").Value = '''1.424'
$ws.Range("E29").Value = '  -0.97%  '

$ws.Range("D30").Value = '''1.483'
$ws.Range("E30").Value = '  +1.12%  '

$ws.Range("D31").Value = '''4.082'
$ws.Range("E31").Value = '  -0.48%  '

$ws.Range("D32").Value = '''4.082'
$ws.Range("E32").Value = '  +0.37%  '

$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").Value = '''1.141'
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("D35").Value = '''0.6948'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").Value = '''2.581'
$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("D37").Value = '''1.262.88'
$ws.Range("E37").Value = '  +2.39%  '

$ws.Range("D38").Value = '''2.831'
$ws.Range("E38").Value = '  +3.74%  '

$ws.Range("D39").Value = '''0.01829'
$ws.Range("E39").Value = '  +1.59%  '

$ws.Range("D40").Value = '''6.718'
$ws.Range("E40").Value = '  +5.35%  '

$ws.Range("D41").Value = '''0.9229'
$ws.Range("E41").Value = '  +2.23%  '

$ws.Range("D42").Value = '''0.9995'
$ws.Range("E42").Value = '  +0.06%  '

$ws.Range("D43").Value = '''2.005.07'
$ws.Range("E43").Value = '  +1.28%  '

$ws.Range("D44").Value = '''101.30'
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").Value = '''65.95'
$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("E46").Value = '  +3.92%  '

$ws.Range("D47").Value = '''1.726'
$ws.Range("E47").Value = '  +2.70%  '

$ws.Range("D48").Value = '''7.054'
$ws.Range("E48").Value = '  -1.82%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.1155'
$ws.Range("E49").Value = '  +0.79%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.014'
$ws.Range("E50").Value = '  +0.51%  '

$ws.Range("D51").Value = '''0.3936'
$ws.Range("E51").Value = '  -1.40%  '
